$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): C1, D1, E1
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Row 2: C2 becomes text (same as D2), E2 becomes numeric 1
$ws.Range("C2").Value = "o__Fusobacteriales"
$ws.Range("E2").Value = 1

# Row 3: C3 becomes text (same as D3), E3 becomes numeric 1
$ws.Range("C3").Value = "o__Fusobacteriales"
$ws.Range("E3").Value = 1
